$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row-by-row rewrite of data rows 2-6 (Maladies/Résumés/URL/Date/Pays)
# Text columns (A, B, C, E) are set directly.
# Date column (D) holds plain text dates like "2023-10-25"; Excel would
# normally auto-convert such a literal to a date serial + apply a date
# number format, so we stage it through a scratch cell pre-formatted as
# Text, copy/paste-special the value only (keeps style untouched), then
# clear the scratch cell.

# Row 2
$ws.Range("A2").Value = 'eastern equine encephalitis, west nile virus'
$ws.Range("B2").Value = 'The first EEE positive mosquito pool was detected in Gloucester County this year (https://www.nj.gov/health/cd/statistics/arboviral-stats/). “Vaccinated animals are much less likely to contract deadly diseases such as EEE and West Nile Virus.” For more information about EEE in horses, visit the New Jersey Department of Agriculture website at: http://www.nj.gov/agriculture/divisions/ah/diseases/diseaseworksheets.html EEE and West Nile virus, like other viral diseases affecting a horse’s neurological system, must be reported to the state veterinarian at 609-671-6400 within 48 hours of diagnosis.'
$ws.Range("C2").Value = 'https://lc.cx/RqYUj8'
$ws.Range("Z1").NumberFormat = "@"
$ws.Range("Z1").Value = '2023-10-25'
$ws.Range("Z1").Copy()
$ws.Range("D2").PasteSpecial(-4163)
$ws.Range("Z1").Clear()
$ws.Range("E2").Value = 'New Jersey, Gloucester County'

# Row 3
$ws.Range("A3").Value = 'bluetongue'
$ws.Range("B3").Value = 'But the disease has resurfaced this year and animal health experts reported that a new strain, named BTV-3, had been confirmed on more than 700 Dutch farms by early October - with a new case also reported in Belgium this week. “The existing BTV-8 serotype vaccine will not offer cross-protection against this new BTV-3 strain, making any likely outbreak difficult to control. "Hence why it’s so important that we follow the advice to take action and prioritise good biosecurity measures while remaining extremely vigilant to the disease at this stage." Dr Henry added: "It remains extremely difficult to protect against midges and a vector-borne disease.'
$ws.Range("C3").Value = 'https://lc.cx/nLtrtC'
$ws.Range("Z1").NumberFormat = "@"
$ws.Range("Z1").Value = '2023-10-11'
$ws.Range("Z1").Copy()
$ws.Range("D3").PasteSpecial(-4163)
$ws.Range("Z1").Clear()
$ws.Range("E3").Value = 'Belgium'

# Row 4
$ws.Range("A4").Value = 'eastern equine encephalitis, west nile virus'
$ws.Range("B4").Value = 'The first EEE positive mosquito pool was detected in Gloucester County this year (https://www.nj.gov/health/cd/statistics/arboviral-stats/). “Vaccinated animals are much less likely to contract deadly diseases such as EEE and West Nile Virus.” For more information about EEE in horses, visit the New Jersey Department of Agriculture website at: http://www.nj.gov/agriculture/divisions/ah/diseases/diseaseworksheets.html EEE and West Nile virus, like other viral diseases affecting a horse’s neurological system, must be reported to the state veterinarian at 609-671-6400 within 48 hours of diagnosis.'
$ws.Range("C4").Value = ' https://lc.cx/RqYUj8'
$ws.Range("Z1").NumberFormat = "@"
$ws.Range("Z1").Value = '2023-10-25'
$ws.Range("Z1").Copy()
$ws.Range("D4").PasteSpecial(-4163)
$ws.Range("Z1").Clear()
$ws.Range("E4").Value = 'New Jersey, Gloucester County'

# Row 5
$ws.Range("A5").Value = 'fièvre aphteuse'
$ws.Range("B5").Value = 'Fièvre Aphteuse en Algérie : Une Campagne Préventive Cruciale pour Protéger le Cheptel National L’Algérie est en alerte face à une menace qui plane sur son secteur agricole. Les services vétérinaires sont mobilisés dans tout le pays pour garantir l’efficacité de la campagne. La fièvre aphteuse peut entraîner des pertes économiques considérables si elle n’est pas maîtrisée. Cependant, la riposte préventive par le biais d’une campagne de vaccination nationale démontre l’engagement du pays à protéger son cheptel national et à préserver son industrie agricole.'
$ws.Range("C5").Value = 'https://www.algerie-focus.com/fievre-aphteuse-en-algerie-une-campagne-preventive-cruciale-pour-proteger-le-cheptel-national/'
$ws.Range("Z1").NumberFormat = "@"
$ws.Range("Z1").Value = '2024-04-09'
$ws.Range("Z1").Copy()
$ws.Range("D5").PasteSpecial(-4163)
$ws.Range("Z1").Clear()
$ws.Range("E5").Value = 'Algeria'

# Row 6
$ws.Range("A6").Value = 'السل'
$ws.Range("B6").Value = 'كانت البداية بإعلان صحيفة محلية بيع 10 أبقار صادرة عن مركب فلاحي بسليانة شمال غربي تونس، مما أثار حفيظة عمادة البياطرة التي قالت إن الأبقار المعروضة للبيع مصابة بالسل، ودعا عميد الأطباء البياطرة أحمد رجب وزارة الفلاحة إلى التدخل من أجل إبطال عملية البيع، مما أعاد فتح ملف الأمراض الحيوانية التي تشكل خطراً على الإنسان. وقال رجب لـ"اندبندنت عربية" إن القانون رقم 95 لعام 2005 المتعلق بتربية الماشية والمنتجات الحيوانية في تونس يحظر على كل مالك أو مربي الاتجار بالحيوانات المصابة بأمراض مثل السل والأمراض شديدة العدوى التي تسبب خسائر اقتصادية ويمكن أن تنتقل من الحيوان إلى الإنسان، منوهاً بأن الأبقار المصابة بالسل لا تحمل علامات ويمكنها إطلاق جرثومة المرض في الحليب ومن ثم نقلها إلى البشر. يشار إلى أن تونس تعرف مئات الإصابات بالسل لدى البشر جراء انتقاله من الحيوان، وتصاعدت التحذيرات من أخطار صحية ناجمة عن انتشار السل الحيواني الذي يتسبب فيه استهلاك الألبان غير المعقمة أو الخاضعة للرقابة، وأسهم نقص الحليب المعلب في ارتفاع الإقبال على الحليب الخام الذي قد يهدد صحة الإنسان. يذكر أن مصالح وزارة الفلاحة في كل المحافظات تقدم سنوياً جدولاً زمنياً لتلقيح القطيع الحيواني لدى الخواص لمجابهة أمراض السل والكلب مجاناً.'
$ws.Range("C6").Value = 'https://lc.cx/GFQ2zw'
$ws.Range("Z1").NumberFormat = "@"
$ws.Range("Z1").Value = '2023-11-18'
$ws.Range("Z1").Copy()
$ws.Range("D6").PasteSpecial(-4163)
$ws.Range("Z1").Clear()
$ws.Range("E6").Value = 'Sulaylana, Tunisia, Baton'

